$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update GPS points for site 9 (rows 21-22 => site_ID 9B/9C)
$ws.Range("B21").Value = 40.24724
$ws.Range("C21").Value = -111.67156
$ws.Range("B22").Value = 40.24724
$ws.Range("C22").Value = -111.67156

# Touch M19 (wrap-text formatting) which extends the used range to column M
$ws.Range("M19").WrapText = $true

$ws.Range("M19").Select()
